$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its text formatting (some new values look numeric,
# e.g. "6.58", and would otherwise be auto-converted to a Number by Excel).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.646.55"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "2.590.47"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "508.06"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "154.00"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -6.53%  "
$ws.Range("D9").Value = "2.598.64"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "6.58"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "3.045.35"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "60.603.84"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "21.64"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("D18").Value = "2.599.72"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "346.05"
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "6.15"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "60.02"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "19.37"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "153.66"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "0.854"
$ws.Range("E37").Value = "  +9.80%  "
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("D40").Value = "3.75"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "36.04"
$ws.Range("D42").Value = "296.16"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "0.625"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0557"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "19.87"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "4.87"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").Value = "0.0234"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "1.999.19"
$ws.Range("E51").Value = "  +0.04%  "

# Restore default style on column D (writing while NumberFormat="@" bumps the
# style index even for cells whose format reverts to General).
$dRange.Style = "Normal"
